# Update the "dSF" column (column F) values for the rows that were repulled.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    4  = -3
    6  = -2
    8  = 1
    9  = -2
    10 = -2
    20 = -5
    21 = -4
    23 = -5
    24 = 0
    29 = -3
    35 = -7
    36 = -12
    39 = -3
    45 = 3
    49 = 0
    52 = -3
    54 = -2
    58 = 0
    63 = -4
    65 = -5
    67 = -2
    68 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
